$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Row 30: D30/E30 switch from numeric to the "N/A" text markers ---
# (copy format+value from sibling cells that already use shared strings 20/"0" and 29/"***.*")
$ws.Range("C30").Copy($ws.Range("D30"))
$ws.Range("M30").Copy($ws.Range("E30"))

# --- Numeric cell updates (crime-stat table, rows 14-30) ---
$ws.Range("J14").Value = 35
$ws.Range("K14").Value = -68.571428571428
$ws.Range("L14").Value = -65.625
$ws.Range("N14").Value = -88.043478260869
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 91
$ws.Range("J15").Value = 125
$ws.Range("K15").Value = -27.2
$ws.Range("L15").Value = -22.222222222222
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -61.440677966101
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 121
$ws.Range("G16").Value = 132
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 875
$ws.Range("J16").Value = 975
$ws.Range("K16").Value = -10.256410256410
$ws.Range("L16").Value = 22.720897615708
$ws.Range("M16").Value = -36.084733382030
$ws.Range("N16").Value = -82.387278582930
$ws.Range("C17").Value = 50
$ws.Range("D17").Value = 61
$ws.Range("E17").Value = -18.032786885245
$ws.Range("F17").Value = 227
$ws.Range("G17").Value = 222
$ws.Range("H17").Value = 2.252252252252
$ws.Range("I17").Value = 1831
$ws.Range("J17").Value = 1812
$ws.Range("K17").Value = 1.048565121412
$ws.Range("L17").Value = 13.726708074534
$ws.Range("M17").Value = 64.954954954955
$ws.Range("N17").Value = -24.650205761316
$ws.Range("C18").Value = 21
$ws.Range("D18").Value = 30
$ws.Range("E18").Value = -30
$ws.Range("F18").Value = 92
$ws.Range("G18").Value = 104
$ws.Range("H18").Value = -11.538461538461
$ws.Range("I18").Value = 746
$ws.Range("J18").Value = 746
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 27.086882453151
$ws.Range("M18").Value = -44.535315985130
$ws.Range("N18").Value = -87.510463753557
$ws.Range("C19").Value = 80
$ws.Range("D19").Value = 71
$ws.Range("E19").Value = 12.676056338028
$ws.Range("F19").Value = 313
$ws.Range("G19").Value = 298
$ws.Range("H19").Value = 5.033557046979
$ws.Range("I19").Value = 2338
$ws.Range("J19").Value = 2455
$ws.Range("K19").Value = -4.765784114052
$ws.Range("L19").Value = 41.525423728813
$ws.Range("M19").Value = 29.529085872576
$ws.Range("N19").Value = -56.759755872017
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = 42
$ws.Range("E20").Value = 9.523809523809
$ws.Range("G20").Value = 116
$ws.Range("H20").Value = 41.379310344827
$ws.Range("I20").Value = 1172
$ws.Range("J20").Value = 1020
$ws.Range("K20").Value = 14.901960784313
$ws.Range("L20").Value = 70.101596516690
$ws.Range("M20").Value = 5.585585585585
$ws.Range("N20").Value = -90.503200713070
$ws.Range("C21").Value = 228
$ws.Range("D21").Value = 237
$ws.Range("E21").Value = -3.797468354430
$ws.Range("F21").Value = 925
$ws.Range("G21").Value = 892
$ws.Range("H21").Value = 3.699551569506
$ws.Range("I21").Value = 7064
$ws.Range("J21").Value = 7168
$ws.Range("K21").Value = -1.450892857142
$ws.Range("L21").Value = 30.814814814814
$ws.Range("M21").Value = 2.764038405586
$ws.Range("N21").Value = -77.536807962603
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 79
$ws.Range("J22").Value = 74
$ws.Range("K22").Value = 6.756756756756
$ws.Range("L22").Value = 23.4375
$ws.Range("M22").Value = -1.25
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = -4.545454545454
$ws.Range("I23").Value = 162
$ws.Range("J23").Value = 151
$ws.Range("K23").Value = 7.284768211920
$ws.Range("L23").Value = 6.578947368421
$ws.Range("M23").Value = 62
$ws.Range("D24").Value = 260
$ws.Range("E24").Value = -33.461538461538
$ws.Range("F24").Value = 687
$ws.Range("G24").Value = 905
$ws.Range("H24").Value = -24.088397790055
$ws.Range("I24").Value = 6076
$ws.Range("J24").Value = 6553
$ws.Range("K24").Value = -7.279108805127
$ws.Range("L24").Value = 40.518038852914
$ws.Range("M24").Value = 51.294820717131
$ws.Range("C25").Value = 92
$ws.Range("D25").Value = 80
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 373
$ws.Range("G25").Value = 325
$ws.Range("H25").Value = 14.769230769230
$ws.Range("I25").Value = 2972
$ws.Range("J25").Value = 2664
$ws.Range("K25").Value = 11.561561561561
$ws.Range("L25").Value = 32.088888888888
$ws.Range("M25").Value = -2.268990463663
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 165
$ws.Range("J26").Value = 199
$ws.Range("K26").Value = -17.085427135678
$ws.Range("L26").Value = -9.340659340659
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 60
$ws.Range("F27").Value = 31
$ws.Range("G27").Value = 31
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 278
$ws.Range("J27").Value = 283
$ws.Range("K27").Value = -1.766784452296
$ws.Range("L27").Value = 6.923076923076
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("I28").Value = 60
$ws.Range("J28").Value = 129
$ws.Range("K28").Value = -53.488372093023
$ws.Range("L28").Value = -60.264900662251
$ws.Range("M28").Value = -53.846153846153
$ws.Range("N28").Value = -81.927710843373
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -50
$ws.Range("I29").Value = 44
$ws.Range("J29").Value = 98
$ws.Range("K29").Value = -55.102040816326
$ws.Range("L29").Value = -64.227642276422
$ws.Range("M29").Value = -58.095238095238
$ws.Range("N29").Value = -85.382059800664
$ws.Range("G30").Value = 4
